$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 239
$ws1.Range("F5").Value = 1868
$ws1.Range("F7").Value = 639
$ws1.Range("F8").Value = 11
$ws1.Range("F9").Value = 124
$ws1.Range("F10").Value = 141
$ws1.Range("F11").Value = 622
$ws1.Range("F12").Value = 17
$ws1.Range("F14").Value = 460
$ws1.Range("F17").Value = 176
$ws1.Range("F18").Value = 244

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 23
$ws2.Range("F12").Value = 201

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6307
$ws3.Range("F4").Value = 1958
$ws3.Range("F5").Value = 167

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6307
$ws4.Range("F4").Value = 1958
$ws4.Range("F6").Value = 167
$ws4.Range("F12").Value = 239
$ws4.Range("F16").Value = 1868
$ws4.Range("F20").Value = 23
$ws4.Range("F21").Value = 639
$ws4.Range("F22").Value = 11
$ws4.Range("F23").Value = 124
$ws4.Range("F24").Value = 201
$ws4.Range("F25").Value = 141
$ws4.Range("F26").Value = 622
$ws4.Range("F27").Value = 17
$ws4.Range("F30").Value = 460
$ws4.Range("F35").Value = 176
$ws4.Range("F41").Value = 244
